$wb = $excel.ActiveWorkbook

# Sheet1: ETNotification "NewAdd" row
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("K2").Value = "30-05-2024"
$ws1.Range("N2").Value = "30-05-2024 12:51:09 PM"
$ws1.Range("AG2").Value = "ET458"

# AK2 needs to stay a text value ("3") rather than become numeric -
# force text formatting, assign, then drop back to General so no
# stray number-format-only diff is left behind.
$ak2 = $ws1.Range("AK2")
$ak2.NumberFormat = "@"
$ak2.Value = "3"
$ak2.NumberFormat = "General"

# Sheet2: ETNotification "SummaryAdd" row
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("K2").Value = "30-05-2024"
$ws2.Range("N2").Value = "30-05-2024 12:57:11 PM"
$ws2.Range("AG2").Value = "ET459"

# Sheet3: ETNotification "DuplicateAdd" row
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("K2").Value = "30-05-2024"
$ws3.Range("N2").Value = "30-05-2024 12:57:11 PM"
$ws3.Range("AG2").Value = "ET460"

# Sheet4: ETNotification "EditSave" row
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("K2").Value = "30-05-2024"
$ws4.Range("N2").Value = "30-05-2024 12:57:11 PM"
$ws4.Range("AG2").Value = "ET460"
